$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("animal_clinic")
$ws.Range("G14").Value = "11 - Oogontsteking "
